$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.233.80'
$ws.Range('E2').Value = '  +3.69%  '
$ws.Range('D3').Value = '1.606.01'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '''212.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.63%  '
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  +2.21%  '
$ws.Range('D8').Value = '''0.249'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.69%  '
$ws.Range('D9').Value = '''0.0617'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('D10').Value = '''18.04'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('D11').Value = '''0.0822'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.20%  '
$ws.Range('D12').Value = '1.829.98'
$ws.Range('E12').Value = '  +2.64%  '
$ws.Range('D13').Value = '1.607.86'
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = '''4.01'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '''0.511'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '26.202.99'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('D17').Value = '''60.71'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').Value = '0.0₃0724'
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('E19').Value = '  +12.63%  '
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = '''4.25'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').Value = '''9.33'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('D23').Value = '''5.99'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('E24').Value = '  +11.75%  '
$ws.Range('D25').Value = '''141.90'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.52%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('E27').Value = '  -3.65%  '
$ws.Range('D28').Value = '''15.25'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.79%  '
$ws.Range('D29').Value = '''6.45'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('D31').Value = '''0.0471'
$ws.Range('D31').Style = "Normal"
$ws.Range('E32').Value = '  +3.34%  '
$ws.Range('D33').Value = '''3.00'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('D34').Value = '''1.47'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('E35').Value = '  +2.52%  '
$ws.Range('D36').Value = '1.109.60'
$ws.Range('E36').Value = '  +2.42%  '
$ws.Range('D37').Value = '''0.0162'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +7.68%  '
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').Value = '''2.34'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('D40').Value = '''0.782'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('D41').Value = '''0.495'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.14%  '
$ws.Range('E42').Value = '  +1.81%  '
$ws.Range('D43').Value = '1.744.10'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('D44').Value = '''92.71'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('D45').Value = '''5.09'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '0.0₆0106'
$ws.Range('E46').Value = '  -5.45%  '
$ws.Range('D47').Value = '''1.51'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +6.99%  '
$ws.Range('D48').Value = '''53.55'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('E51').Value = '  +0.00%  '
